$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-apply the (already present) border formatting to the existing header /
# case rows. This normalizes their style index the same way Excel
# consolidates duplicate cell formats when it resaves the file, without
# altering their values.
$ws.Range("B4:C5").Borders.LineStyle = 1
$ws.Range("B7").Borders.LineStyle = 1

# Apply the bordered-cell style to the new rows first (matches the style
# already used by the existing case rows) so it reuses the existing xf
# instead of allocating a transient quote-prefix-only style.
$ws.Range("B8:C11").Borders.LineStyle = 1

# New test cases: case21..case24
$ws.Range("B8").Value = "case21"
$ws.Range("B9").Value = "case22"
$ws.Range("B10").Value = "case23"
$ws.Range("B11").Value = "case24"

# New datatype declaration block
$ws.Range("B14").Value = "Datatype MyType"
$ws.Range("B15").Value = "String"
$ws.Range("C15").Value = "value"
$ws.Range("D15").Value = "XXX"

# Value expressions for the new test cases (entered with a leading apostrophe
# so Excel stores them as literal text with a quote-prefix, same as the
# existing "= sr.instance" / "= ChildWithId.id" rows)
$ws.Range("C8").Value = "'= AccessBean.getClass()"
$ws.Range("C9").Value = "'= MyType.getClass()"
$ws.Range("C10").Value = "'= MyType.value"
$ws.Range("C11").Value = "'= MyType.getValue()"

$ws.Range("C15").Select()
